$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 201 }

$ws.Range("C2:C$lastRow").Value = 45205
